$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 17:55"

# Refresh country statistics and re-sort order (rows that moved / got new data)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 190256
$ws.Range("C4").Value = 1726
$ws.Range("D4").Value = 7274
$ws.Range("E4").Value = 178869
$ws.Range("F4").Value = 4576
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 4113

$ws.Range("A20").Value = "Brasil"
$ws.Range("B20").Value = 5908
$ws.Range("C20").Value = 191
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 5577
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 204

$ws.Range("A25").Value = "Chequia"
$ws.Range("B25").Value = 3508
$ws.Range("C25").Value = 200
$ws.Range("D25").Value = 48
$ws.Range("E25").Value = 3421
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 39

$ws.Range("A33").Value = "Luxemburgo"
$ws.Range("B33").Value = 2319
$ws.Range("C33").Value = 141
$ws.Range("D33").Value = 80
$ws.Range("E33").Value = 2210
$ws.Range("F33").Value = 31
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 29

$ws.Range("A34").Value = "Filipinas"
$ws.Range("B34").Value = 2311
$ws.Range("C34").Value = 227
$ws.Range("D34").Value = 50
$ws.Range("E34").Value = 2165
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 96

$ws.Range("A35").Value = "Ecuador"
$ws.Range("B35").Value = 2302
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 58
$ws.Range("E35").Value = 2165
$ws.Range("F35").Value = 100
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 79

$ws.Range("A71").Value = "Libano"
$ws.Range("B71").Value = 479
$ws.Range("C71").Value = 16
$ws.Range("D71").Value = 43
$ws.Range("E71").Value = 422
$ws.Range("F71").Value = 5
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 14

$ws.Range("A72").Value = "Bosnia y Herzegovina"
$ws.Range("B72").Value = 457
$ws.Range("C72").Value = 37
$ws.Range("D72").Value = 19
$ws.Range("E72").Value = 425
$ws.Range("F72").Value = 1
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 13

$ws.Range("A74").Value = "Moldavia"
$ws.Range("B74").Value = 423
$ws.Range("C74").Value = 70
$ws.Range("D74").Value = 22
$ws.Range("E74").Value = 396
$ws.Range("F74").Value = 44
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 5

$ws.Range("A77").Value = "Tunez"
$ws.Range("B77").Value = 394
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 5
$ws.Range("E77").Value = 379
$ws.Range("F77").Value = 10
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 10

$ws.Range("A79").Value = "Kazajistan"
$ws.Range("B79").Value = 380
$ws.Range("C79").Value = 37
$ws.Range("D79").Value = 26
$ws.Range("E79").Value = 351
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 3

$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 354
$ws.Range("C81").Value = 25
$ws.Range("D81").Value = 17
$ws.Range("E81").Value = 326
$ws.Range("F81").Value = 4
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 11

$ws.Range("A87").Value = "Reunion"
$ws.Range("B87").Value = 281
$ws.Range("C87").Value = 34
$ws.Range("D87").Value = 40
$ws.Range("E87").Value = 241
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0

$ws.Range("A88").Value = "Jordania"
$ws.Range("B88").Value = 274
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 30
$ws.Range("E88").Value = 239
$ws.Range("F88").Value = 5
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 5

$ws.Range("A89").Value = "Burkina Faso"
$ws.Range("B89").Value = 261
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 32
$ws.Range("E89").Value = 215
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 14

$ws.Range("A90").Value = "Albania"
$ws.Range("B90").Value = 259
$ws.Range("C90").Value = 16
$ws.Range("D90").Value = 67
$ws.Range("E90").Value = 177
$ws.Range("F90").Value = 7
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 15

$ws.Range("A94").Value = "Cuba"
$ws.Range("B94").Value = 212
$ws.Range("C94").Value = 26
$ws.Range("D94").Value = 12
$ws.Range("E94").Value = 194
$ws.Range("F94").Value = 3
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 6

$ws.Range("A95").Value = "Oman"
$ws.Range("B95").Value = 210
$ws.Range("C95").Value = 18
$ws.Range("D95").Value = 34
$ws.Range("E95").Value = 175
$ws.Range("F95").Value = 3
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 1

$ws.Range("A96").Value = "Afganistan"
$ws.Range("B96").Value = 196
$ws.Range("C96").Value = 22
$ws.Range("D96").Value = 5
$ws.Range("E96").Value = 187
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 4

$ws.Range("A97").Value = "Ghana"
$ws.Range("B97").Value = 195
$ws.Range("C97").Value = 34
$ws.Range("D97").Value = 31
$ws.Range("E97").Value = 159
$ws.Range("F97").Value = 1
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 5

$ws.Range("A98").Value = "Senegal"
$ws.Range("B98").Value = 190
$ws.Range("C98").Value = 15
$ws.Range("D98").Value = 45
$ws.Range("E98").Value = 144
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 1

$ws.Range("A99").Value = "Malta"
$ws.Range("B99").Value = 188
$ws.Range("C99").Value = 19
$ws.Range("D99").Value = 2
$ws.Range("E99").Value = 186
$ws.Range("F99").Value = 2
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

$ws.Range("A113").Value = "Georgia"
$ws.Range("B113").Value = 117
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = 23
$ws.Range("E113").Value = 94
$ws.Range("F113").Value = 6
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A114").Value = "Bolivia"
$ws.Range("B114").Value = 115
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = 1
$ws.Range("E114").Value = 107
$ws.Range("F114").Value = 3
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 0
